$d = $word.ActiveDocument

# 0. Header/footer picture-anchor runs: fix the east-Asian language tag
#    from zh-TW to es-MX (affects the (empty) runs that host the logo
#    drawings in the header and footer).
$sec = $d.Sections.Item(1)
$sec.Headers.Item(1).Range.LanguageIDFarEast = "es-MX"
$sec.Footers.Item(1).Range.LanguageIDFarEast = "es-MX"

# 1. Title paragraph: expand the intro sentence with the contract date.
$d.Content.Find.Execute(
    "CONTRATO INDIVIDUAL DE TRABAJO QUE CELEBRARAN, POR UNA ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "CONTRATO INDIVIDUAL DE TRABAJO QUE CELEBRA EL DIA 8 DE ENERO DEL 2019,, POR UNA ",
    2) | Out-Null

# 2. "SEPTIMA" clause: consolidate/expand the work start & end dates.
$d.Content.Find.Execute(
    "EL TRABAJO EMPEZARÁ EL 08 DE ENERO HASTA EL 21 DE AGOSTO. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "EL TRABAJO EMPEZARÁ EL 08 DE ENERO DEL 2019 HASTA EL 21 DE AGOSTO DEL 2019.",
    2) | Out-Null

# 3. Signature line: insert "DEL" before the year.
$d.Content.Find.Execute(
    "LA CIUDAD DE OCOSINGO, EL DIA 08 DE ENERO 2019.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "LA CIUDAD DE OCOSINGO, EL DIA 08 DE ENERO DEL 2019.",
    2) | Out-Null
